# Regenerate merged AHB files
#
# 1. Rename the header row (row 1) of the "AHB-Diff" sheet:
#      *_old -> *_FV2310   (columns A-J)
#      *_new -> *_FV2404   (columns L-U)
#      "diff" (column K) stays as-is
# 2. Freeze the header row (split/freeze at row 1).
# 3. Turn the A1:U56 range into an Excel Table ("Table1") with an
#    autofilter, picking up the new header names as its column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10): "<name>_old" -> "<name>_FV2310"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2310"
}

# Column K (11): "diff" unchanged
$ws.Cells.Item(1, 11).Value = "diff"

# Columns L-U (12-21): "<name>_new" -> "<name>_FV2404"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2404"
}

# Freeze the top row (pane split after row 1).
$ws.Activate()
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into a table with an autofilter / header row.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U56"), $null, 1)
$tbl.Name = "Table1"
